# Auto-generated edit script applying the Phantom_Profits market-data refresh diff.
# For each affected (sheet, row), updates columns H-N (currentAveragePrice..LeveProfitHQ)
# to the new scraped values. Cells that the diff removes are cleared; cells the diff
# newly introduces are written for the first time.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3523
$ws.Range("I4").Value = 2284.5
$ws.Range("K4").Value = 2284.5
$ws.Range("M4").Value = -2170.5

$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H41").Value = 3080
$ws.Range("J41").Value = 2950
$ws.Range("L41").Value = 2950
$ws.Range("N41").Value = -3830

$ws.Range("H47").Value = 9990
$ws.Range("I47").Value = 9990
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 9990
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -9018
$ws.Range("N47").ClearContents()

$ws.Range("H64").Value = 4615.3335
$ws.Range("J64").Value = 3498.5
$ws.Range("L64").Value = 3498.5
$ws.Range("N64").Value = -3994.5

$ws.Range("H67").Value = 4615.3335
$ws.Range("J67").Value = 3498.5
$ws.Range("L67").Value = 3498.5
$ws.Range("N67").Value = -5214.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10442.036
$ws.Range("I32").Value = 8976.963
$ws.Range("K32").Value = 8976.963
$ws.Range("M32").Value = -8689.963

$ws.Range("H74").Value = 5310991.5
$ws.Range("I74").Value = 6902789.5
$ws.Range("K74").Value = 6902789.5
$ws.Range("M74").Value = -6901915.5

$ws.Range("H77").Value = 5310991.5
$ws.Range("I77").Value = 6902789.5
$ws.Range("K77").Value = 34513947.5
$ws.Range("M77").Value = -34509579.5

$ws.Range("H96").Value = 19661.572
$ws.Range("J96").Value = 19661.572
$ws.Range("L96").Value = 19661.572
$ws.Range("N96").Value = -25153.572

$ws.Range("H113").Value = 12000
$ws.Range("J113").Value = 12000
$ws.Range("L113").Value = 12000
$ws.Range("N113").Value = -20678

$ws.Range("H132").Value = 6219.6313
$ws.Range("I132").Value = 6219.6313
$ws.Range("K132").Value = 18658.8939
$ws.Range("M132").Value = -16128.8939

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 611.25
$ws.Range("I22").Value = 570.55554
$ws.Range("J22").Value = 733.3333
$ws.Range("K22").Value = 570.55554
$ws.Range("L22").Value = 733.3333
$ws.Range("M22").Value = -397.55554
$ws.Range("N22").Value = -1079.3333

$ws.Range("H30").Value = 5000
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 5000
$ws.Range("N30").Value = -5250

$ws.Range("H54").Value = 4997.5
$ws.Range("I54").Value = 4997.5
$ws.Range("K54").Value = 4997.5
$ws.Range("M54").Value = -4513.5

$ws.Range("H134").Value = 1113
$ws.Range("I134").Value = 965.1667
$ws.Range("K134").Value = 2895.5001
$ws.Range("M134").Value = -360.5001000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 610.6667
$ws.Range("I16").Value = 543.5
$ws.Range("K16").Value = 543.5
$ws.Range("M16").Value = -256.5

$ws.Range("H31").Value = 2873.75
$ws.Range("I31").Value = 2165
$ws.Range("K31").Value = 2165
$ws.Range("M31").Value = -1870

$ws.Range("H34").Value = 2873.75
$ws.Range("I34").Value = 2165
$ws.Range("K34").Value = 2165
$ws.Range("M34").Value = -1963

$ws.Range("H86").Value = 5002
$ws.Range("I86").Value = 5002
$ws.Range("K86").Value = 5002
$ws.Range("M86").Value = -3879

$ws.Range("H89").Value = 5002
$ws.Range("I89").Value = 5002
$ws.Range("K89").Value = 25010
$ws.Range("M89").Value = -19394

$ws.Range("H99").Value = 12124
$ws.Range("I99").Value = 3228
$ws.Range("K99").Value = 3228
$ws.Range("M99").Value = -1730

$ws.Range("H105").Value = 1699.5
$ws.Range("I105").Value = 1699.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1699.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 47.5
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 610.6667
$ws.Range("I113").Value = 543.5
$ws.Range("K113").Value = 543.5
$ws.Range("M113").Value = 1626.5

$ws.Range("H122").Value = 4473.9
$ws.Range("J122").Value = 2998.5
$ws.Range("L122").Value = 8995.5
$ws.Range("N122").Value = -13895.5

$ws.Range("H126").Value = 12124
$ws.Range("I126").Value = 3228
$ws.Range("K126").Value = 9684
$ws.Range("M126").Value = -7214

$ws.Range("H141").Value = 602371.3
$ws.Range("J141").Value = 602371.3
$ws.Range("L141").Value = 602371.3
$ws.Range("N141").Value = -612731.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.782608
$ws.Range("I2").Value = 43.555557
$ws.Range("J2").Value = 50.5
$ws.Range("K2").Value = 261.333342
$ws.Range("L2").Value = 303
$ws.Range("M2").Value = -148.333342
$ws.Range("N2").Value = -529

$ws.Range("H7").Value = 102243.5
$ws.Range("J7").Value = 497
$ws.Range("L7").Value = 1491
$ws.Range("N7").Value = -1715

$ws.Range("H34").Value = 129873.875
$ws.Range("I34").Value = 998.5
$ws.Range("J34").Value = 172832.33
$ws.Range("K34").Value = 2995.5
$ws.Range("L34").Value = 518496.99
$ws.Range("M34").Value = -2911.5
$ws.Range("N34").Value = -518664.99

$ws.Range("H55").Value = 499.5
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H131").Value = 1422.75
$ws.Range("I131").Value = 1494
$ws.Range("J131").Value = 1399
$ws.Range("K131").Value = 4482
$ws.Range("L131").Value = 4197
$ws.Range("M131").Value = 558
$ws.Range("N131").Value = -14277

$ws.Range("H137").Value = 7224.5
$ws.Range("J137").Value = 9082.833000000001
$ws.Range("L137").Value = 27248.499
$ws.Range("N137").Value = -37448.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1766.2222
$ws.Range("I107").Value = 737.125
$ws.Range("K107").Value = 737.125
$ws.Range("M107").Value = 1182.875

$ws.Range("H132").Value = 47623228
$ws.Range("I132").Value = 5295.2
$ws.Range("J132").Value = 166668060
$ws.Range("K132").Value = 15885.6
$ws.Range("L132").Value = 500004180
$ws.Range("M132").Value = -13355.6
$ws.Range("N132").Value = -500009240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2381
$ws.Range("I82").Value = 1823.75
$ws.Range("J82").Value = 3124
$ws.Range("K82").Value = 1823.75
$ws.Range("L82").Value = 3124
$ws.Range("M82").Value = -1462.75
$ws.Range("N82").Value = -3846

$ws.Range("H85").Value = 2381
$ws.Range("I85").Value = 1823.75
$ws.Range("J85").Value = 3124
$ws.Range("K85").Value = 1823.75
$ws.Range("L85").Value = 3124
$ws.Range("M85").Value = -575.75
$ws.Range("N85").Value = -5620

$ws.Range("H122").Value = 3518.2
$ws.Range("I122").Value = 3147.75
$ws.Range("K122").Value = 9443.25
$ws.Range("M122").Value = -6993.25

$ws.Range("H132").Value = 27777.5
$ws.Range("I132").Value = 5555
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 16665
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -14135
$ws.Range("N132").Value = -155060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5258.857
$ws.Range("I136").Value = 5616.1035
$ws.Range("K136").Value = 16848.3105
$ws.Range("M136").Value = -14298.3105
